# Refresh the cryptocurrency price/volume table with the latest feed values
# (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these Price cells to remain text (they look numeric) so the
# stored value matches the source feed's literal formatting exactly.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D10", "D11", "D13", "D15", "D16", "D17", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D34", "D35", "D36", "D38", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the refreshed crypto feed.
$ws.Range("D2").Value = '29.006.73'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '1.830.12'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("D4").Value = '0.9980'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '243.78'
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").Value = '0.6320'
$ws.Range("E6").Value = '  +1.37%  '
$ws.Range("D7").Value = '0.9993'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '0.07515'
$ws.Range("E8").Value = '  -0.60%  '
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("D10").Value = '22.89'
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("D11").Value = '0.07726'
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("D12").Value = '1.840.61'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").Value = '4.993'
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("D15").Value = '83.09'
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("D16").Value = '0.000009677'
$ws.Range("E16").Value = '  +6.81%  '
$ws.Range("D17").Value = '6.080'
$ws.Range("E17").Value = '  +1.25%  '
$ws.Range("D18").Value = '29.041.48'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '12.55'
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("D20").Value = '226.40'
$ws.Range("E20").Value = '  +0.68%  '
$ws.Range("D21").Value = '0.9986'
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").Value = '0.9994'
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = '159.79'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").Value = '0.1407'
$ws.Range("E25").Value = '  +3.57%  '
$ws.Range("D26").Value = '8.538'
$ws.Range("E26").Value = '  +1.74%  '
$ws.Range("D27").Value = '17.91'
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("D28").Value = '1.498'
$ws.Range("E28").Value = '  +0.20%  '
$ws.Range("D29").Value = '4.120'
$ws.Range("E29").Value = '  +1.63%  '
$ws.Range("D30").Value = '4.074'
$ws.Range("E30").Value = '  +1.21%  '
$ws.Range("D31").Value = '1.198'
$ws.Range("E31").Value = '  -0.50%  '
$ws.Range("E32").Value = '  +3.00%  '
$ws.Range("E33").Value = '  +1.09%  '
$ws.Range("D34").Value = '0.7436'
$ws.Range("E34").Value = '  +1.55%  '
$ws.Range("D35").Value = '1.138'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("D36").Value = '2.653'
$ws.Range("E36").Value = '  +1.65%  '
$ws.Range("D37").Value = '1.246.07'
$ws.Range("E37").Value = '  -2.58%  '
$ws.Range("D38").Value = '2.753'
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("E39").Value = '  +0.33%  '
$ws.Range("D40").Value = '6.588'
$ws.Range("E40").Value = '  +3.55%  '
$ws.Range("D41").Value = '0.9011'
$ws.Range("E41").Value = '  +1.08%  '
$ws.Range("D42").Value = '0.9989'
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("D43").Value = '101.54'
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("D44").Value = '1.986.21'
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.00000000123'
$ws.Range("E45").Value = '  +3.27%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '64.88'
$ws.Range("E46").Value = '  +2.23%  '
$ws.Range("D47").Value = '0.5102'
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("D48").Value = '0.4069'
$ws.Range("E48").Value = '  +2.63%  '
$ws.Range("D49").Value = '8.995'
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("D50").Value = '6.764'
$ws.Range("E50").Value = '  +1.37%  '
$ws.Range("E51").Value = '  +0.25%  '
